# 2017-02-13 snapshot - chunk 30
# Update the STEO report month references (January 2017 -> February 2017)
# and refresh the underlying cooling-degree-day data for Fig29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / source text updates -------------------------------------------------
$ws.Range("A2").Value2  = "Short-Term Energy Outlook, February 2017"
$ws.Range("A34").Value2 = "Source: Short-Term Energy Outlook, February 2017."

# --- Data table updates (B27:F32) -------------------------------------------------
$ws.Range("B27").Value2 = 53.356033363999998
$ws.Range("C27").Value2 = 42.665723094999997
$ws.Range("D27").Value2 = 42.886724031
$ws.Range("E27").Value2 = 39.891996532999997
$ws.Range("F27").Value2 = 39.470820000000003

$ws.Range("B28").Value2 = 125.9657717
$ws.Range("C28").Value2 = 97.464444354999998
$ws.Range("D28").Value2 = 126.32850126
$ws.Range("E28").Value2 = 126.74356014999999
$ws.Range("F28").Value2 = 115.6249

$ws.Range("B29").Value2 = 255.16354622
$ws.Range("C29").Value2 = 270.31827881999999
$ws.Range("D29").Value2 = 247.66019671999999
$ws.Range("E29").Value2 = 248.12415813999999
$ws.Range("F29").Value2 = 250.42339999999999

$ws.Range("B30").Value2 = 336.01668622
$ws.Range("C30").Value2 = 383.58931858
$ws.Range("D30").Value2 = 357.12762593999997
$ws.Range("E30").Value2 = 357.56511886999999
$ws.Range("F30").Value2 = 346.50330000000002

$ws.Range("B31").Value2 = 315.52240819000002
$ws.Range("C31").Value2 = 361.71782325999999
$ws.Range("D31").Value2 = 331.88751832999998
$ws.Range("E31").Value2 = 332.36211229000003
$ws.Range("F31").Value2 = 323.428

$ws.Range("B32").Value2 = 223.37339832999999
$ws.Range("C32").Value2 = 220.16345195
$ws.Range("D32").Value2 = 182.98378316
$ws.Range("E32").Value2 = 183.49430365000001
$ws.Range("F32").Value2 = 187.47909999999999
